$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (row 9) onto the two
# new rows before filling them in, so the new cells pick up the same
# per-column styles (A -> s1, B -> s0, C/D -> s2, E/F/G -> s0).
$ws.Range("A9:G9").Copy()
$ws.Range("A12:G12").PasteSpecial(-4122)
$ws.Range("A9:G9").Copy()
$ws.Range("A13:G13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 12: new review row
$ws.Range("A12").Value = "com.hamxa.shaynachim"
$ws.Range("B12").Value = "bitcoin"
$ws.Range("C12").Value = "stavsade45@gmail.com"
$ws.Range("D12").Value = "galiatia942@gmail.com"
$ws.Range("E12").Value = "27/5/2019 15:59"
$ws.Range("F12").Value = "Just loved this bitcoin guide for beginners app, really useful and give some idea about this hard topic"
$ws.Range("G12").Value = "no"

# Row 13: new review row
$ws.Range("A13").Value = "com.singleton.strechy"
$ws.Range("B13").Value = "taxi"
$ws.Range("C13").Value = "itamaramir2@gmail.com"
$ws.Range("D13").Value = "cohenyossi408@gmail.com"
$ws.Range("E13").Value = "27/5/2019 15:59"
$ws.Range("F13").Value = "Cool Game Ever !! I had play this in long time ago like 2 years ago but since that i like this game so much!! Should try this game its so cool much better than another car game.READ ALL REVIEWS"
$ws.Range("G13").Value = "confirm"

$ws.Rows.Item(12).RowHeight = 13.8
$ws.Rows.Item(13).RowHeight = 13.8

$ws.Range("A1").Select()
$ws.Range("D18").Select()
